$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting existing rows 50-63 down to 51-64
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new record
$ws.Cells.Item(50, 1).Value = 7
$ws.Cells.Item(50, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(50, 3).Value = "Ñuble"
$ws.Cells.Item(50, 4).Value = 45218
$ws.Cells.Item(50, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(50, 5).Value = 16
$ws.Cells.Item(50, 6).Value = 300000000
$ws.Cells.Item(50, 7).Value = "Espárragos"
$ws.Cells.Item(50, 8).Value = "Sin especificar"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 300
$ws.Cells.Item(50, 11).Value = 1300
$ws.Cells.Item(50, 12).Value = 1300
$ws.Cells.Item(50, 13).Value = 1300
$ws.Cells.Item(50, 14).Value = "$/kilo"
$ws.Cells.Item(50, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(50, 16).Value = 1300
$ws.Cells.Item(50, 17).Value = 1
$ws.Cells.Item(50, 18).Value = "Hortaliza"
